$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New oil-drop measurement rows appended below the existing dataset
# (A = charge value, B = uncertainty value, D = A*10/1.626 computed charge)
$newRows = @(
    @(15, 0.08826, 0.081244),
    @(16, 0.14076, 0.12266),
    @(17, 0.38213, 0.34023),
    @(18, 0.31022, 0.28239),
    @(19, 0.23096, 0.20501),
    @(20, 0.12446, 0.0717),
    @(21, 0.34867, 0.30979),
    @(22, 0.34917, 0.31994),
    @(23, 0.26865, 0.28305),
    @(24, 0.36111, 0.32122),
    @(25, 0.19791, 0.17281)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 4).Formula = "=A$r*10/1.626"
}

# Match the author's final selection (cell B15)
$ws.Range("B15").Select() | Out-Null
